$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Val)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "39.511.89"
Set-TextValue $ws "E2" "  +1.81%  "
Set-TextValue $ws "D3" "2.159.61"
Set-TextValue $ws "E3" "  +2.56%  "
Set-TextValue $ws "D4" "0.999"
Set-TextValue $ws "E4" "  -0.10%  "
Set-TextValue $ws "D5" "227.83"
Set-TextValue $ws "E5" "  +0.08%  "
Set-TextValue $ws "D6" "0.629"
Set-TextValue $ws "E6" "  +2.13%  "
Set-TextValue $ws "D7" "63.25"
Set-TextValue $ws "E7" "  +1.49%  "
Set-TextValue $ws "E8" "  +0.02%  "
Set-TextValue $ws "D9" "0.393"
Set-TextValue $ws "E9" "  +0.84%  "
Set-TextValue $ws "D10" "0.0848"
Set-TextValue $ws "E10" "  +0.63%  "
Set-TextValue $ws "E11" "  +0.01%  "
Set-TextValue $ws "D12" "15.98"
Set-TextValue $ws "E12" "  +1.10%  "
Set-TextValue $ws "D13" "2.475.87"
Set-TextValue $ws "E13" "  +2.43%  "
Set-TextValue $ws "D14" "21.95"
Set-TextValue $ws "E14" "  -0.20%  "
Set-TextValue $ws "D15" "0.808"
Set-TextValue $ws "E15" "  -0.30%  "
Set-TextValue $ws "D16" "5.48"
Set-TextValue $ws "E16" "  -0.86%  "
Set-TextValue $ws "D17" "2.148.54"
Set-TextValue $ws "E17" "  +1.94%  "
Set-TextValue $ws "D18" "39.472.21"
Set-TextValue $ws "E18" "  +1.82%  "
Set-TextValue $ws "D19" "71.93"
Set-TextValue $ws "E19" "  +0.44%  "
Set-TextValue $ws "D20" "6.11"
Set-TextValue $ws "E20" "  -0.41%  "
Set-TextValue $ws "D21" "0.0₃0844"
Set-TextValue $ws "E21" "  -0.10%  "
Set-TextValue $ws "D22" "227.89"
Set-TextValue $ws "E22" "  -0.22%  "
Set-TextValue $ws "E23" "  +0.09%  "
Set-TextValue $ws "B24" "Toncoin"
Set-TextValue $ws "C24" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D24" "2.37"
Set-TextValue $ws "E24" "  +0.24%  "
Set-TextValue $ws "B25" "PancakeSwap"
Set-TextValue $ws "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D25" "2.38"
Set-TextValue $ws "E25" "  +2.81%  "
Set-TextValue $ws "B26" "Cosmos"
Set-TextValue $ws "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D26" "9.62"
Set-TextValue $ws "E26" "  -0.50%  "
Set-TextValue $ws "B27" "Monero"
Set-TextValue $ws "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D27" "171.64"
Set-TextValue $ws "E27" "  -0.34%  "
Set-TextValue $ws "E28" "  +0.75%  "
Set-TextValue $ws "D29" "19.72"
Set-TextValue $ws "E29" "  +1.91%  "
Set-TextValue $ws "D30" "1.42"
Set-TextValue $ws "E30" "  -0.66%  "
Set-TextValue $ws "D31" "2.69"
Set-TextValue $ws "E31" "  +5.00%  "
Set-TextValue $ws "E32" "  +1.61%  "
Set-TextValue $ws "D33" "4.61"
Set-TextValue $ws "E33" "  +0.97%  "
Set-TextValue $ws "E34" "  -1.29%  "
Set-TextValue $ws "D35" "6.95"
Set-TextValue $ws "E35" "  -3.46%  "
Set-TextValue $ws "D36" "0.0619"
Set-TextValue $ws "E36" "  +0.08%  "
Set-TextValue $ws "D37" "2.41"
Set-TextValue $ws "E37" "  +0.94%  "
Set-TextValue $ws "D38" "3.60"
Set-TextValue $ws "E38" "  +1.53%  "
Set-TextValue $ws "D39" "1.00"
Set-TextValue $ws "E39" "  +0.09%  "
Set-TextValue $ws "D40" "4.71"
Set-TextValue $ws "E40" "  +12.95%  "
Set-TextValue $ws "D41" "101.81"
Set-TextValue $ws "E41" "  -0.75%  "
Set-TextValue $ws "D42" "0.0226"
Set-TextValue $ws "E42" "  -0.61%  "
Set-TextValue $ws "D43" "17.66"
Set-TextValue $ws "E43" "  -2.79%  "
Set-TextValue $ws "D44" "1.514.65"
Set-TextValue $ws "E44" "  -0.90%  "
Set-TextValue $ws "D45" "1.21"
Set-TextValue $ws "E45" "  +1.25%  "
Set-TextValue $ws "D46" "0.0924"
Set-TextValue $ws "E46" "  +0.85%  "
Set-TextValue $ws "E47" "  -0.09%  "
Set-TextValue $ws "B48" "ARBITRUM"
Set-TextValue $ws "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D48" "1.10"
Set-TextValue $ws "E48" "  +1.61%  "
Set-TextValue $ws "B49" "FraxShare"
Set-TextValue $ws "C49" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D49" "7.75"
Set-TextValue $ws "E49" "  -0.71%  "
Set-TextValue $ws "D50" "0.000192"
Set-TextValue $ws "E50" "  +35.42%  "
Set-TextValue $ws "E51" "  +0.56%  "

Write-Host "Applied 104 cell updates"
